$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("Baptiste", "raph", "yass", "vincent")
$startRow = 19
$startId = 18

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $id = $startId + $i
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).NumberFormat = "General"
}
